$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing cell contents (keeps per-cell number formatting / styles intact)
# so the new "apoio_medio" / "contribuicoes" / "media_contribuicoes" columns can be
# written alongside the existing ones exactly like the expanded analysis output.
$ws.Cells.ClearContents()

# ---- Header row (row 1) ----
$headers = @(
  "geral_modalidade",
  "mencoes_ficcao_cientifica",
  "total",
  "total_sucesso",
  "particip",
  "taxa_sucesso",
  "arrecadado_sucesso",
  "media_sucesso",
  "std_sucesso",
  "min_sucesso",
  "max_sucesso",
  "apoio_medio",
  "contribuicoes",
  "media_contribuicoes"
)
for ($col = 1; $col -le $headers.Length; $col++) {
  $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# The three new header cells (L1:N1) need the same bold / bordered / centered
# formatting that the rest of the header row (A1:K1) already has. Copy that
# formatting from the existing K1 header cell.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Data rows (rows 2-7) ----
# Columns: A geral_modalidade, B mencoes_ficcao_cientifica, C total, D total_sucesso,
# E particip, F taxa_sucesso, G arrecadado_sucesso, H media_sucesso, I std_sucesso,
# J min_sucesso, K max_sucesso, L apoio_medio, M contribuicoes, N media_contribuicoes
$data = @(
  @("aon",  $false, 1039, 651,  77.82771535580524, 62.65640038498557, 18661633.68066395, 28666.10396415353, 44525.51900988264, 41.81688448509265, 679297.6600721752, 92.37335010030465, 202024, 310.3287250384025),
  @("aon",  $true,  296,  179,  22.17228464419476, 60.47297297297297, 5401646.146659081, 30176.79411541386, 46623.43714083682, 54.53892516702949, 537544.5528256212, 87.79024763378376, 61529,  343.7374301675978),
  @("flex", $false, 1116, 1045, 76.02179836512262, 93.63799283154121, 13073618.32949033, 12510.63954975151, 30558.59191293883, 23.05352861032933, 475290.9541363961, 89.64725872903666, 145834, 139.5540669856459),
  @("flex", $true,  352,  338,  23.97820163487738, 96.02272727272727, 5288513.608068768, 15646.4899647005,  42686.00050640347, 10.77163914429046, 708972.7845446636, 91.47778329877478, 57812,  171.0414201183432),
  @("sub",  $false, 613,  135,  89.61988304093568, 22.02283849918434, 41147.50408794444, 304.7963265773662, 685.6746025059557, 1.087396962410123, 5087.076865717208,  19.42752789799076, 2118,   15.68888888888889),
  @("sub",  $true,  71,   17,   10.38011695906433, 23.94366197183098, 2039.453666840368, 119.9678627553158, 157.2645751483118, 2.022084306600051, 538.4389998789497, 22.66059629822632, 90,     5.294117647058823)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $rowNum = $i + 2
  $rowVals = $data[$i]
  for ($col = 1; $col -le $rowVals.Length; $col++) {
    $ws.Cells.Item($rowNum, $col).Value = $rowVals[$col - 1]
  }
}
